$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "20.570.81"
$ws.Range("E2").Value = "  +0.94%  "

$ws.Range("D3").Value = "1.476.14"

$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9513"
$ws.Range("E5").Value = "  +5.55%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "279.31"
$ws.Range("E6").Value = "  -0.12%  "

$ws.Range("E7").Value = "  -1.86%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3067"
$ws.Range("E8").Value = "  -3.40%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "39.92"
$ws.Range("E9").Value = "  +1.07%  "

$ws.Range("E10").Value = "  +1.51%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06670"
$ws.Range("E11").Value = "  +1.31%  "

$ws.Range("E12").Value = "  -0.16%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.525"
$ws.Range("E13").Value = "  +0.25%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.11"
$ws.Range("E14").Value = "  +1.97%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.222"
$ws.Range("E15").Value = "  +0.79%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9537"
$ws.Range("E16").Value = "  +6.05%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001033"
$ws.Range("E17").Value = "  +0.42%  "

$ws.Range("D18").Value = "1.474.98"
$ws.Range("E18").Value = "  -0.16%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.05943"
$ws.Range("E19").Value = "  +5.10%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.81"
$ws.Range("E20").Value = "  +0.17%  "

$ws.Range("E21").Value = "  -2.69%  "

$ws.Range("E22").Value = "  -0.57%  "

$ws.Range("E23").Value = "  -0.49%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.275"
$ws.Range("E24").Value = "  -0.40%  "

$ws.Range("D25").Value = "20.625.47"
$ws.Range("E25").Value = "  +0.38%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "143.22"
$ws.Range("E26").Value = "  +4.60%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.121"
$ws.Range("E27").Value = "  -5.50%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.25"
$ws.Range("E28").Value = "  -0.44%  "

$ws.Range("D29").Value = "1.636.21"
$ws.Range("E29").Value = "  -0.62%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "113.81"
$ws.Range("E30").Value = "  +1.18%  "

$ws.Range("E31").Value = "  +0.63%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.020"
$ws.Range("E32").Value = "  +0.91%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.8121"
$ws.Range("E33").Value = "  -1.35%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07969"
$ws.Range("E34").Value = "  +2.30%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.517"
$ws.Range("E35").Value = "  +2.23%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.225"
$ws.Range("E36").Value = "  +5.90%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05846"
$ws.Range("E37").Value = "  -3.14%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.740"
$ws.Range("E38").Value = "  -1.70%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02056"
$ws.Range("E39").Value = "  +1.27%  "

$ws.Range("E40").Value = "  +0.28%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9543"
$ws.Range("E41").Value = "  +3.98%  "

$ws.Range("E42").Value = "  +1.38%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.430"
$ws.Range("E43").Value = "  +7.83%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5312"
$ws.Range("E44").Value = "  -0.13%  "

$ws.Range("B45").Value = "PancakeSwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.539"
$ws.Range("E45").Value = "  -0.81%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.27"
$ws.Range("E46").Value = "  +0.46%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "118.17"
$ws.Range("E47").Value = "  -3.18%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5198"
$ws.Range("E48").Value = "  -0.69%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.819"
$ws.Range("E49").Value = "  +0.44%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06479"
$ws.Range("E50").Value = "  +0.99%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9872"
$ws.Range("E51").Value = "  -0.67%  "
